$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace column A (company names) and column B (values)
$names = @(
    "Mettler Toledo",
    "Leggett & Platt",
    "Zoetis",
    "PPG Industries",
    "Vulcan Materials",
    "Jacobs Engineering Group",
    "Procter & Gamble",
    "Fluor Corp.",
    "Anadarko Petroleum Corp",
    "Public Storage"
)
$values = @(396, 410, 461, 396, 769, 759, 509, 392, 780, 565)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Update the chart: style and series references now cover row 10 too
$chart = $ws.ChartObjects().Item(1).Chart
$chart.ChartStyle = 1
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(,'Sheet1'!`$A`$1:`$A`$10,'Sheet1'!`$B`$1:`$B`$10,1)"
